$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date-serial value that must be bumped
# from 46060 (2026-02-07) to 46061 (2026-02-08) for every data row (2-96).
$ws.Range("C2:C96").Value = 46061
